# Update "想去人数" (F column) figures across sheets, matching the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 828
$ws1.Range("F6").Value  = 3718
$ws1.Range("F7").Value  = 2498
$ws1.Range("F9").Value  = 2327
$ws1.Range("F14").Value = 627
$ws1.Range("F17").Value = 20
$ws1.Range("F18").Value = 43
$ws1.Range("F19").Value = 258
$ws1.Range("F21").Value = 419
$ws1.Range("F24").Value = 462
$ws1.Range("F25").Value = 654
$ws1.Range("F26").Value = 72
$ws1.Range("F28").Value = 347
$ws1.Range("F30").Value = 1602
$ws1.Range("F31").Value = 774
$ws1.Range("F32").Value = 796
$ws1.Range("F33").Value = 1893
$ws1.Range("F35").Value = 493
$ws1.Range("F36").Value = 72
$ws1.Range("F37").Value = 544
$ws1.Range("F38").Value = 1181
$ws1.Range("F40").Value = 396

# --- 演出 (sheet 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 57
$ws2.Range("F8").Value = 5

# --- 全部类型 (sheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 828
$ws4.Range("F6").Value  = 3719
$ws4.Range("F7").Value  = 2498
$ws4.Range("F9").Value  = 2327
$ws4.Range("F14").Value = 627
$ws4.Range("F17").Value = 20
$ws4.Range("F18").Value = 43
$ws4.Range("F19").Value = 258
$ws4.Range("F21").Value = 419
$ws4.Range("F24").Value = 462
$ws4.Range("F25").Value = 654
$ws4.Range("F26").Value = 72
$ws4.Range("F27").Value = 57
$ws4.Range("F31").Value = 347
$ws4.Range("F33").Value = 1602
$ws4.Range("F34").Value = 774
$ws4.Range("F36").Value = 796
$ws4.Range("F37").Value = 1893
$ws4.Range("F41").Value = 5
$ws4.Range("F42").Value = 493
$ws4.Range("F43").Value = 72
$ws4.Range("F44").Value = 544
$ws4.Range("F45").Value = 1181
$ws4.Range("F47").Value = 396
